$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stats")

$ws.Cells.Item(3, 4).Value = 565870383.2374853
$ws.Cells.Item(5, 4).Value = 329.6370627282254
$ws.Cells.Item(7, 4).Value = 2448249180.756474
$ws.Cells.Item(9, 4).Value = 1426.181140907407
$ws.Cells.Item(11, 4).Value = 354879520.7513012
$ws.Cells.Item(13, 4).Value = 206.7283362200003
$ws.Cells.Item(15, 4).Value = 2803128701.507775
$ws.Cells.Item(17, 4).Value = 1632.909477127407
$ws.Cells.Item(19, 4).Value = 1385185508.822446
$ws.Cells.Item(21, 4).Value = 806.9135547429833
$ws.Cells.Item(23, 4).Value = 1951055892.059932
$ws.Cells.Item(25, 4).Value = 1136.550617471209
$ws.Cells.Item(27, 4).Value = 826591921.680119
$ws.Cells.Item(29, 4).Value = 481.5154516103398
$ws.Cells.Item(31, 4).Value = 25480887.76772453
$ws.Cells.Item(33, 4).Value = 14.84340804585843
$ws.Cells.Item(35, 4).Value = 852072809.4478436
$ws.Cells.Item(37, 4).Value = 496.3588596561983
$ws.Cells.Item(39, 4).Value = 2055226995.564199
$ws.Cells.Item(41, 4).Value = 1197.233518710614
$ws.Cells.Item(43, 4).Value = 4197439.824399191
$ws.Cells.Item(45, 4).Value = 2.445138985322621
$ws.Cells.Item(47, 4).Value = 350682080.9269021
$ws.Cells.Item(49, 4).Value = 204.2831972346777
$ws.Cells.Item(51, 4).Value = 1360965249.830945
$ws.Cells.Item(53, 4).Value = 784.8985354769284
$ws.Cells.Item(55, 4).Value = 4541721265.441251
$ws.Cells.Item(57, 4).Value = 2619.31035361268
$ws.Cells.Item(59, 4).Value = 759583831.5953854
$ws.Cells.Item(61, 4).Value = 438.068670059893
$ws.Cells.Item(63, 4).Value = 5301305097.036636
$ws.Cells.Item(65, 4).Value = 3057.379023672572
$ws.Cells.Item(67, 4).Value = 2196795589.333132
$ws.Cells.Item(69, 4).Value = 1266.94024040947
$ws.Cells.Item(71, 4).Value = 3557760839.164076
$ws.Cells.Item(73, 4).Value = 2051.838775886398
$ws.Cells.Item(75, 4).Value = 1703831566.423668
$ws.Cells.Item(77, 4).Value = 982.6370668548799
$ws.Cells.Item(79, 4).Value = 39712691.44889126
$ws.Cells.Item(81, 4).Value = 22.90318093129402
$ws.Cells.Item(83, 4).Value = 1743544257.87256
$ws.Cells.Item(85, 4).Value = 1005.540247786174
$ws.Cells.Item(87, 4).Value = 4507647282.917536
$ws.Cells.Item(89, 4).Value = 2599.659139899432
$ws.Cells.Item(91, 4).Value = 44120664.28010463
$ws.Cells.Item(93, 4).Value = 25.44535562684281
$ws.Cells.Item(95, 4).Value = 715463167.3152808
$ws.Cells.Item(97, 4).Value = 412.6233144330502
$ws.Cells.Item(99, 4).Value = 1698587797.396736
$ws.Cells.Item(101, 4).Value = 988.1247554271557
$ws.Cells.Item(103, 4).Value = 6254818799.136025
$ws.Cells.Item(105, 4).Value = 3638.635168349724
$ws.Cells.Item(107, 4).Value = 793301879.8173603
$ws.Cells.Item(109, 4).Value = 461.4899666510098
$ws.Cells.Item(111, 4).Value = 7048120678.953385
$ws.Cells.Item(113, 4).Value = 4100.125135000733
$ws.Cells.Item(115, 4).Value = 2665471924.885136
$ws.Cells.Item(117, 4).Value = 1550.593262186198
$ws.Cells.Item(119, 4).Value = 4364059722.281873
$ws.Cells.Item(121, 4).Value = 2538.718017613354
$ws.Cells.Item(123, 4).Value = 2635901521.320776
$ws.Cells.Item(125, 4).Value = 1533.391179470959
$ws.Cells.Item(127, 4).Value = 48159435.35073736
$ws.Cells.Item(129, 4).Value = 28.015937916421
$ws.Cells.Item(131, 4).Value = 2684060956.671514
$ws.Cells.Item(133, 4).Value = 1561.40711738738
$ws.Cells.Item(135, 4).Value = 7089011641.785209
$ws.Cells.Item(137, 4).Value = 4123.91276182185
$ws.Cells.Item(139, 4).Value = 65922821.53550412
$ws.Cells.Item(141, 4).Value = 38.34948773720844
$ws.Cells.Item(143, 4).Value = 727379058.2818561
$ws.Cells.Item(145, 4).Value = 423.1404789138012
$ws.Cells.Item(147, 4).Value = 2367594408.420507
$ws.Cells.Item(149, 4).Value = 1365.826030152444
$ws.Cells.Item(151, 4).Value = 8873890289.731684
$ws.Cells.Item(153, 4).Value = 5119.200443845567
$ws.Cells.Item(155, 4).Value = 1000559328.409302
$ws.Cells.Item(157, 4).Value = 577.2061171427437
$ws.Cells.Item(159, 4).Value = 9874449618.140985
$ws.Cells.Item(161, 4).Value = 5696.40656098831
$ws.Cells.Item(163, 4).Value = 3386336814.594242
$ws.Cells.Item(165, 4).Value = 1953.521663924649
$ws.Cells.Item(167, 4).Value = 5753931223.01475
$ws.Cells.Item(169, 4).Value = 3319.347694077093
$ws.Cells.Item(171, 4).Value = 4061409553.121212
$ws.Cells.Item(173, 4).Value = 2342.959954219288
$ws.Cells.Item(175, 4).Value = 59108842.00502516
$ws.Cells.Item(177, 4).Value = 34.09891269192957
$ws.Cells.Item(179, 4).Value = 4120518395.126236
$ws.Cells.Item(181, 4).Value = 2377.058866911218
$ws.Cells.Item(183, 4).Value = 10646852463.23067
$ws.Cells.Item(185, 4).Value = 6141.992978930179
$ws.Cells.Item(187, 4).Value = 41452819.32045235
$ws.Cells.Item(189, 4).Value = 23.913445415193
$ws.Cells.Item(191, 4).Value = 959106509.0888498
$ws.Cells.Item(193, 4).Value = 553.2926717275507
$ws.Cells.Item(195, 4).Value = 3465750206.07974
$ws.Cells.Item(197, 4).Value = 2008.059134209944
$ws.Cells.Item(199, 4).Value = 11155608716.82801
$ws.Cells.Item(201, 4).Value = 6463.570843103962
$ws.Cells.Item(203, 4).Value = 1395361978.369195
$ws.Cells.Item(205, 4).Value = 808.4741252494794
$ws.Cells.Item(207, 4).Value = 12550970695.1972
$ws.Cells.Item(209, 4).Value = 7272.044968353442
$ws.Cells.Item(211, 4).Value = 4026392247.819052
$ws.Cells.Item(213, 4).Value = 2332.89569368322
$ws.Cells.Item(215, 4).Value = 7492142453.898792
$ws.Cells.Item(217, 4).Value = 4340.954827893164
$ws.Cells.Item(219, 4).Value = 4986661817.249826
$ws.Cells.Item(221, 4).Value = 2889.276842219737
$ws.Cells.Item(223, 4).Value = 72166424.0485834
$ws.Cells.Item(225, 4).Value = 41.81329824054028
$ws.Cells.Item(227, 4).Value = 5058828241.298409
$ws.Cells.Item(229, 4).Value = 2931.090140460277
$ws.Cells.Item(231, 4).Value = 12787019102.64856
$ws.Cells.Item(233, 4).Value = 7408.811651614946
$ws.Cells.Item(235, 4).Value = 36587528.69525979
$ws.Cells.Item(237, 4).Value = 21.19885070360852
$ws.Cells.Item(239, 4).Value = 1358774449.673935
$ws.Cells.Item(241, 4).Value = 787.2752745458708
$ws.Cells.Item(243, 4).Value = 4887820646.135211
$ws.Cells.Item(245, 4).Value = 2831.82192688992
$ws.Cells.Item(247, 4).Value = 15035449047.22926
$ws.Cells.Item(249, 4).Value = 8710.981309481163
$ws.Cells.Item(251, 4).Value = 2037983561.439964
$ws.Cells.Item(253, 4).Value = 1180.732059080396
$ws.Cells.Item(255, 4).Value = 17073432608.66922
$ws.Cells.Item(257, 4).Value = 9891.713368561561
$ws.Cells.Item(259, 4).Value = 4720591516.152766
$ws.Cells.Item(261, 4).Value = 2734.935573771899
$ws.Cells.Item(263, 4).Value = 9608412162.287977
$ws.Cells.Item(265, 4).Value = 5566.757500661819
$ws.Cells.Item(267, 4).Value = 7380649242.960949
$ws.Cells.Item(269, 4).Value = 4276.074323108895
$ws.Cells.Item(271, 4).Value = 84371203.42029516
$ws.Cells.Item(273, 4).Value = 48.88154479084631
$ws.Cells.Item(275, 4).Value = 7465020446.381245
$ws.Cells.Item(277, 4).Value = 4324.955867899742
$ws.Cells.Item(279, 4).Value = 19994460239.89693
$ws.Cells.Item(281, 4).Value = 11584.04839761024
$ws.Cells.Item(283, 4).Value = 176821989.1494206
$ws.Cells.Item(285, 4).Value = 102.4440997902707
$ws.Cells.Item(287, 4).Value = 1861161572.290543
$ws.Cells.Item(289, 4).Value = 1078.287959290125
$ws.Cells.Item(291, 4).Value = 6661808773.102343
$ws.Cells.Item(293, 4).Value = 3863.93934514537
$ws.Cells.Item(295, 4).Value = 20418128382.66599
$ws.Cells.Item(297, 4).Value = 11842.79109459819
$ws.Cells.Item(299, 4).Value = 2594149026.014614
$ws.Cells.Item(301, 4).Value = 1504.641581616653
$ws.Cells.Item(303, 4).Value = 23012277408.6806
$ws.Cells.Item(305, 4).Value = 13347.43267621484
$ws.Cells.Item(307, 4).Value = 5797930692.11448
$ws.Cells.Item(309, 4).Value = 3362.8783540201
$ws.Cells.Item(311, 4).Value = 12459739465.21682
$ws.Cells.Item(313, 4).Value = 7226.817699165469
$ws.Cells.Item(315, 4).Value = 10448869313.23742
$ws.Cells.Item(317, 4).Value = 6060.485766975632
$ws.Cells.Item(319, 4).Value = 103668630.2263575
$ws.Cells.Item(321, 4).Value = 60.12921007373918
$ws.Cells.Item(323, 4).Value = 10552537943.46378
$ws.Cells.Item(325, 4).Value = 6120.614977049371
$ws.Cells.Item(327, 4).Value = 27934739811.17214
$ws.Cells.Item(329, 4).Value = 16202.52756107272
$ws.Cells.Item(331, 4).Value = 259351516.0167249
$ws.Cells.Item(333, 4).Value = 150.4273930837324
$ws.Cells.Item(335, 4).Value = 2334797509.997889
$ws.Cells.Item(337, 4).Value = 1354.214188532921
$ws.Cells.Item(339, 4).Value = 10139192537.30381
$ws.Cells.Item(341, 4).Value = 5874.529498492469
$ws.Cells.Item(343, 4).Value = 26562199530.05717
$ws.Cells.Item(345, 4).Value = 15389.82755382782
$ws.Cells.Item(347, 4).Value = 4106269335.422499
$ws.Cells.Item(349, 4).Value = 2379.124397820055
$ws.Cells.Item(351, 4).Value = 30668468865.47967
$ws.Cells.Item(353, 4).Value = 17768.95195164787
$ws.Cells.Item(355, 4).Value = 7098879698.986835
$ws.Cells.Item(357, 4).Value = 4113.007820348279
$ws.Cells.Item(359, 4).Value = 17238072236.29065
$ws.Cells.Item(361, 4).Value = 9987.537318840748
$ws.Cells.Item(363, 4).Value = 13299198870.51287
$ws.Cells.Item(365, 4).Value = 7705.400186820078
$ws.Cells.Item(367, 4).Value = 131197758.6761527
$ws.Cells.Item(369, 4).Value = 76.01444598704742
$ws.Cells.Item(371, 4).Value = 13430396629.18902
$ws.Cells.Item(373, 4).Value = 7781.414632807125
$ws.Cells.Item(375, 4).Value = 36954227969.53468
$ws.Cells.Item(377, 4).Value = 21410.8472151348
$ws.Cells.Item(379, 4).Value = 359124291.2995687
$ws.Cells.Item(381, 4).Value = 208.0724115951663
$ws.Cells.Item(383, 4).Value = 3747145044.12293
$ws.Cells.Item(385, 4).Value = 2171.051986224889
$ws.Cells.Item(387, 4).Value = 14019771650.52578
$ws.Cells.Item(389, 4).Value = 8143.253796809007
$ws.Cells.Item(391, 4).Value = 36166308503.07888
$ws.Cells.Item(393, 4).Value = 21006.86347649735
$ws.Cells.Item(395, 4).Value = 5711713231.345447
$ws.Cells.Item(397, 4).Value = 3317.595437133508
$ws.Cells.Item(399, 4).Value = 41878021734.42433
$ws.Cells.Item(401, 4).Value = 24324.45891363086
$ws.Cells.Item(403, 4).Value = 9126918008.739893
$ws.Cells.Item(405, 4).Value = 5301.285326216303
$ws.Cells.Item(407, 4).Value = 23146689659.26567
$ws.Cells.Item(409, 4).Value = 13444.53912302531
$ws.Cells.Item(411, 4).Value = 18562913871.93915
$ws.Cells.Item(413, 4).Value = 10782.09564574741
$ws.Cells.Item(415, 4).Value = 168418203.219502
$ws.Cells.Item(417, 4).Value = 97.8241448581316
$ws.Cells.Item(419, 4).Value = 18731332075.15865
$ws.Cells.Item(421, 4).Value = 10879.91979060555
$ws.Cells.Item(423, 4).Value = 52530605376.93568
$ws.Cells.Item(425, 4).Value = 30511.91291467037
$ws.Cells.Item(427, 4).Value = 569738852.0981594
$ws.Cells.Item(429, 4).Value = 330.9275062524964
$ws.Cells.Item(431, 4).Value = 5141974379.247288
$ws.Cells.Item(433, 4).Value = 2986.667930881013
$ws.Cells.Item(435, 4).Value = 30068075373.95647
$ws.Cells.Item(437, 4).Value = 17357.2388737968
$ws.Cells.Item(439, 4).Value = 67558311597.62813
$ws.Cells.Item(441, 4).Value = 38999.02929357762
$ws.Cells.Item(443, 4).Value = 16953386417.27632
$ws.Cells.Item(445, 4).Value = 9786.591729091004
$ws.Cells.Item(447, 4).Value = 84511698014.90445
$ws.Cells.Item(449, 4).Value = 48785.62102266862
$ws.Cells.Item(451, 4).Value = 18360789863.54364
$ws.Cells.Item(453, 4).Value = 10599.03607429261
$ws.Cells.Item(455, 4).Value = 48428865237.50011
$ws.Cells.Item(457, 4).Value = 27956.27494808941
$ws.Cells.Item(459, 4).Value = 35754365149.90895
$ws.Cells.Item(461, 4).Value = 20639.73330416275
$ws.Cells.Item(463, 4).Value = 328467627.4953835
$ws.Cells.Item(465, 4).Value = 189.6127704164552
$ws.Cells.Item(467, 4).Value = 36082832777.40434
$ws.Cells.Item(469, 4).Value = 20829.34607457921
$ws.Cells.Item(471, 4).Value = 114503741833.9084
$ws.Cells.Item(473, 4).Value = 66098.96956278615
$ws.Cells.Item(475, 4).Value = 3500070761.157085
$ws.Cells.Item(477, 4).Value = 2020.467340228082
$ws.Cells.Item(479, 4).Value = 13453315656.11924
$ws.Cells.Item(481, 4).Value = 7766.124388862921
$ws.Cells.Item(482, 4).Value = 0.008177423883535651
$ws.Cells.Item(483, 4).Value = 0.009401867225265063
$ws.Cells.Item(484, 4).Value = 0.001224443341729412
$ws.Cells.Item(485, 4).Value = 0.00274993990915392
$ws.Cells.Item(486, 4).Value = 0.00008488439352583798
$ws.Cells.Item(487, 4).Value = 0.002834824302679758
$ws.Cells.Item(488, 4).Value = 0.001959521597754708
$ws.Cells.Item(489, 4).Value = 0.004607521324830596
$ws.Cells.Item(490, 4).Value = 0.006567042922585305
$ws.Cells.Item(491, 4).Value = 0.006864853711230775
$ws.Cells.Item(492, 4).Value = -0.00001340022540617041
$ws.Cells.Item(493, 4).Value = -0.001211043116323242
$ws.Cells.Item(494, 4).Value = 0.01044158969010639
$ws.Cells.Item(495, 4).Value = 0.01219588537373859
$ws.Cells.Item(496, 4).Value = 0.001754295683632206
$ws.Cells.Item(497, 4).Value = 0.003889661598052252
$ws.Cells.Item(498, 4).Value = 0.00009185603889536787
$ws.Cells.Item(499, 4).Value = 0.00398151763694762
$ws.Cells.Item(500, 4).Value = 0.003133977023911183
$ws.Cells.Item(501, 4).Value = 0.005080390712879791
$ws.Cells.Item(502, 4).Value = 0.008214367736790974
$ws.Cells.Item(503, 4).Value = 0.01029037668303315
$ws.Cells.Item(504, 4).Value = -0.00009982882552748909
$ws.Cells.Item(505, 4).Value = -0.001654466858104716
$ws.Cells.Item(506, 4).Value = 0.0119525728863239
$ws.Cells.Item(507, 4).Value = 0.01347379531651955
$ws.Cells.Item(508, 4).Value = 0.001521222430195656
$ws.Cells.Item(509, 4).Value = 0.005040713959485842
$ws.Cells.Item(510, 4).Value = 0.00009201763483256039
$ws.Cells.Item(511, 4).Value = 0.005132731594318402
$ws.Cells.Item(512, 4).Value = 0.003248265865405061
$ws.Cells.Item(513, 4).Value = 0.00509279785679609
$ws.Cells.Item(514, 4).Value = 0.008341063722201151
$ws.Cells.Item(515, 4).Value = 0.01356483407561823
$ws.Cells.Item(516, 4).Value = -0.0001298120030550227
$ws.Cells.Item(517, 4).Value = -0.001391410427140633
$ws.Cells.Item(518, 4).Value = 0.01418441817510217
$ws.Cells.Item(519, 4).Value = 0.01578264223063037
$ws.Cells.Item(520, 4).Value = 0.0015982240555282
$ws.Cells.Item(521, 4).Value = 0.006485434610173781
$ws.Cells.Item(522, 4).Value = 0.00009456967356012611
$ws.Cells.Item(523, 4).Value = 0.006580004283733908
$ws.Cells.Item(524, 4).Value = 0.003784824225945553
$ws.Cells.Item(525, 4).Value = 0.005417813720950912
$ws.Cells.Item(526, 4).Value = 0.009202637946896465
$ws.Cells.Item(527, 4).Value = 0.01699810197932478
$ws.Cells.Item(528, 4).Value = -0.00006662344126720546
$ws.Cells.Item(529, 4).Value = -0.001531600614260995
$ws.Cells.Item(530, 4).Value = 0.01529927616649398
$ws.Cells.Item(531, 4).Value = 0.01721261494302797
$ws.Cells.Item(532, 4).Value = 0.001913338776533992
$ws.Cells.Item(533, 4).Value = 0.00685119392922873
$ws.Cells.Item(534, 4).Value = 0.00009883299909468833
$ws.Cells.Item(535, 4).Value = 0.006950026928323418
$ws.Cells.Item(536, 4).Value = 0.004748701110381668
$ws.Cells.Item(537, 4).Value = 0.005513886904322886
$ws.Cells.Item(538, 4).Value = 0.01026258801470455
$ws.Cells.Item(539, 4).Value = 0.01756690909363961
$ws.Cells.Item(540, 4).Value = -0.00004949174202050649
$ws.Cells.Item(541, 4).Value = -0.001863847034513485
$ws.Cells.Item(542, 4).Value = 0.01750221994897444
$ws.Cells.Item(543, 4).Value = 0.01987852306405156
$ws.Cells.Item(544, 4).Value = 0.002376303115077124
$ws.Cells.Item(545, 4).Value = 0.008587823671318152
$ws.Cells.Item(546, 4).Value = 0.00009849346052731308
$ws.Cells.Item(547, 4).Value = 0.008686317131845464
$ws.Cells.Item(548, 4).Value = 0.005684233672642073
$ws.Cells.Item(549, 4).Value = 0.005507972259564023
$ws.Cells.Item(550, 4).Value = 0.0111922059322061
$ws.Cells.Item(551, 4).Value = 0.02327675087240247
$ws.Cells.Item(552, 4).Value = -0.0002071762213551619
$ws.Cells.Item(553, 4).Value = -0.002169126893721963
$ws.Cells.Item(554, 4).Value = 0.01978660622207844
$ws.Cells.Item(555, 4).Value = 0.02230129813125927
$ws.Cells.Item(556, 4).Value = 0.002514691909180832
$ws.Cells.Item(557, 4).Value = 0.01010678843024938
$ws.Cells.Item(558, 4).Value = 0.0001004721153783398
$ws.Cells.Item(559, 4).Value = 0.01020726054562772
$ws.Cells.Item(560, 4).Value = 0.006470482891773825
$ws.Cells.Item(561, 4).Value = 0.005623554693857723
$ws.Cells.Item(562, 4).Value = 0.01209403758563155
$ws.Cells.Item(563, 4).Value = 0.02696793458825862
$ws.Cells.Item(564, 4).Value = -0.0002391321880696474
$ws.Cells.Item(565, 4).Value = -0.002275559721111184
$ws.Cells.Item(566, 4).Value = 0.02125704139546012
$ws.Cells.Item(567, 4).Value = 0.02456842249759844
$ws.Cells.Item(568, 4).Value = 0.003311381102138317
$ws.Cells.Item(569, 4).Value = 0.01065474979436273
$ws.Cells.Item(570, 4).Value = 0.0001049008669360914
$ws.Cells.Item(571, 4).Value = 0.01075965066129882
$ws.Cells.Item(572, 4).Value = 0.008134411035061233
$ws.Cells.Item(573, 4).Value = 0.005674360801238384
$ws.Cells.Item(574, 4).Value = 0.01380877183629962
$ws.Cells.Item(575, 4).Value = 0.02962586961016934
$ws.Cells.Item(576, 4).Value = -0.0002908886900934387
$ws.Cells.Item(577, 4).Value = -0.003020492412044879
$ws.Cells.Item(578, 4).Value = 0.0222759300793821
$ws.Cells.Item(579, 4).Value = 0.02580059637582478
$ws.Cells.Item(580, 4).Value = 0.003524666296442674
$ws.Cells.Item(581, 4).Value = 0.01142841134218476
$ws.Cells.Item(582, 4).Value = 0.0001038489890445706
$ws.Cells.Item(583, 4).Value = 0.01153226033122934
$ws.Cells.Item(584, 4).Value = 0.00864093771679917
$ws.Cells.Item(585, 4).Value = 0.005627398327796272
$ws.Cells.Item(586, 4).Value = 0.01426833604459544
$ws.Cells.Item(587, 4).Value = 0.03232345087273381
$ws.Cells.Item(588, 4).Value = -0.0003503264795785051
$ws.Cells.Item(589, 4).Value = -0.003174339816864169
$ws.Cells.Item(590, 4).Value = 0.02154695761421068
$ws.Cells.Item(591, 4).Value = 0.02621348441030994
$ws.Cells.Item(592, 4).Value = 0.004666526796099271
$ws.Cells.Item(593, 4).Value = 0.01101141277189435
$ws.Cells.Item(594, 4).Value = 0.0001043884044200006
$ws.Cells.Item(595, 4).Value = 0.01111580117631435
$ws.Cells.Item(596, 4).Value = 0.00922016108526918
$ws.Cells.Item(597, 4).Value = 0.00587752214872642
$ws.Cells.Item(598, 4).Value = 0.0150976832339956
$ws.Cells.Item(599, 4).Value = 0.03377289190427531
$ws.Cells.Item(600, 4).Value = -0.0007568511033676168
$ws.Cells.Item(601, 4).Value = -0.003909675692731654
$ws.Cells.Item(602, 4).Value = 0.008177423883535651
$ws.Cells.Item(603, 4).Value = 0.009401867225265063
$ws.Cells.Item(604, 4).Value = 0.008745914545863426
$ws.Cells.Item(605, 4).Value = 0.01657747005333931
$ws.Cells.Item(606, 4).Value = 0.00008488439352583798
$ws.Cells.Item(607, 4).Value = 0.002834824302679758
$ws.Cells.Item(608, 4).Value = 0.01453870303370674
$ws.Cells.Item(609, 4).Value = 0.004607521324830596
$ws.Cells.Item(610, 4).Value = 0.006567042922585305
$ws.Cells.Item(611, 4).Value = 0.04138341584834741
$ws.Cells.Item(612, 4).Value = -0.002566070518277388
$ws.Cells.Item(613, 4).Value = -0.008985354511740578
$ws.Cells.Item(614, 4).Value = 0.01044158969010639
$ws.Cells.Item(615, 4).Value = 0.01219588537373859
$ws.Cells.Item(616, 4).Value = 0.007523891308042997
$ws.Cells.Item(617, 4).Value = 0.01439093679556755
$ws.Cells.Item(618, 4).Value = 0.00009185603889536787
$ws.Cells.Item(619, 4).Value = 0.00398151763694762
$ws.Cells.Item(620, 4).Value = 0.01467007853520475
$ws.Cells.Item(621, 4).Value = 0.005080390712879791
$ws.Cells.Item(622, 4).Value = 0.008214367736790974
$ws.Cells.Item(623, 4).Value = 0.03807224785885416
$ws.Cells.Item(624, 4).Value = -0.00397261207887673
$ws.Cells.Item(625, 4).Value = -0.007744523510258348
$ws.Cells.Item(626, 4).Value = 0.0119525728863239
$ws.Cells.Item(627, 4).Value = 0.01347379531651955
$ws.Cells.Item(628, 4).Value = 0.00447279648167258
$ws.Cells.Item(629, 4).Value = 0.01432330540044378
$ws.Cells.Item(630, 4).Value = 0.00009201763483256039
$ws.Cells.Item(631, 4).Value = 0.005132731594318402
$ws.Cells.Item(632, 4).Value = 0.009957947187551127
$ws.Cells.Item(633, 4).Value = 0.00509279785679609
$ws.Cells.Item(634, 4).Value = 0.008341063722201151
$ws.Cells.Item(635, 4).Value = 0.03854479003034809
$ws.Cells.Item(636, 4).Value = -0.004113264114378734
$ws.Cells.Item(637, 4).Value = -0.004265534941964051
$ws.Cells.Item(638, 4).Value = 0.01418441817510217
$ws.Cells.Item(639, 4).Value = 0.01578264223063037
$ws.Cells.Item(640, 4).Value = 0.004295159852850211
$ws.Cells.Item(641, 4).Value = 0.01353308432356151
$ws.Cells.Item(642, 4).Value = 0.00009456967356012611
$ws.Cells.Item(643, 4).Value = 0.006580004283733908
$ws.Cells.Item(644, 4).Value = 0.01055159569763077
$ws.Cells.Item(645, 4).Value = 0.005417813720950912
$ws.Cells.Item(646, 4).Value = 0.009202637946896465
$ws.Cells.Item(647, 4).Value = 0.0354697504876294
$ws.Cells.Item(648, 4).Value = -0.001236131232238536
$ws.Cells.Item(649, 4).Value = -0.00426990250726572
$ws.Cells.Item(650, 4).Value = 0.01529927616649398
$ws.Cells.Item(651, 4).Value = 0.01721261494302797
$ws.Cells.Item(652, 4).Value = 0.003775568056170293
$ws.Cells.Item(653, 4).Value = 0.01257219440648408
$ws.Cells.Item(654, 4).Value = 0.00009883299909468833
$ws.Cells.Item(655, 4).Value = 0.006950026928323418
$ws.Cells.Item(656, 4).Value = 0.009559250330785696
$ws.Cells.Item(657, 4).Value = 0.005513886904322886
$ws.Cells.Item(658, 4).Value = 0.01026258801470455
$ws.Cells.Item(659, 4).Value = 0.03223592829624263
$ws.Cells.Item(660, 4).Value = -0.001435652607940726
$ws.Cells.Item(661, 4).Value = -0.003751969215804143
$ws.Cells.Item(662, 4).Value = 0.01750221994897444
$ws.Cells.Item(663, 4).Value = 0.01987852306405156
$ws.Cells.Item(664, 4).Value = 0.003707302132163152
$ws.Cells.Item(665, 4).Value = 0.01309499226013791
$ws.Cells.Item(666, 4).Value = 0.00009849346052731308
$ws.Cells.Item(667, 4).Value = 0.008686317131845464
$ws.Cells.Item(668, 4).Value = 0.009115487845344206
$ws.Cells.Item(669, 4).Value = 0.005507972259564023
$ws.Cells.Item(670, 4).Value = 0.0111922059322061
$ws.Cells.Item(671, 4).Value = 0.03549314519967121
$ws.Cells.Item(672, 4).Value = -0.001779731505380834
$ws.Cells.Item(673, 4).Value = -0.003478507565566232
$ws.Cells.Item(674, 4).Value = 0.01978660622207844
$ws.Cells.Item(675, 4).Value = 0.02230129813125927
$ws.Cells.Item(676, 4).Value = 0.003718227529928429
$ws.Cells.Item(677, 4).Value = 0.01330177948973999
$ws.Cells.Item(678, 4).Value = 0.0001004721153783398
$ws.Cells.Item(679, 4).Value = 0.01020726054562772
$ws.Cells.Item(680, 4).Value = 0.009728707592271986
$ws.Cells.Item(681, 4).Value = 0.005623554693857723
$ws.Cells.Item(682, 4).Value = 0.01209403758563155
$ws.Cells.Item(683, 4).Value = 0.03549312639345483
$ws.Cells.Item(684, 4).Value = -0.002861341862805458
$ws.Cells.Item(685, 4).Value = -0.003421422404746317
$ws.Cells.Item(686, 4).Value = 0.02125704139546012
$ws.Cells.Item(687, 4).Value = 0.02456842249759844
$ws.Cells.Item(688, 4).Value = 0.004099285278986587
$ws.Cells.Item(689, 4).Value = 0.01268743544451644
$ws.Cells.Item(690, 4).Value = 0.0001049008669360914
$ws.Cells.Item(691, 4).Value = 0.01075965066129882
$ws.Cells.Item(692, 4).Value = 0.01025315620270902
$ws.Cells.Item(693, 4).Value = 0.005674360801238384
$ws.Cells.Item(694, 4).Value = 0.01380877183629962
$ws.Cells.Item(695, 4).Value = 0.03527781650635808
$ws.Cells.Item(696, 4).Value = -0.001994295890157682
$ws.Cells.Item(697, 4).Value = -0.003807230834083409
$ws.Cells.Item(698, 4).Value = 0.0222759300793821
$ws.Cells.Item(699, 4).Value = 0.02580059637582478
$ws.Cells.Item(700, 4).Value = 0.003912402167151844
$ws.Cells.Item(701, 4).Value = 0.01273154260375151
$ws.Cells.Item(702, 4).Value = 0.0001038489890445706
$ws.Cells.Item(703, 4).Value = 0.01153226033122934
$ws.Cells.Item(704, 4).Value = 0.009663914736279599
$ws.Cells.Item(705, 4).Value = 0.005627398327796272
$ws.Cells.Item(706, 4).Value = 0.01426833604459544
$ws.Cells.Item(707, 4).Value = 0.03600915119037078
$ws.Cells.Item(708, 4).Value = -0.001630922183934227
$ws.Cells.Item(709, 4).Value = -0.003550141239244591
$ws.Cells.Item(710, 4).Value = 0.02154695761421068
$ws.Cells.Item(711, 4).Value = 0.02621348441030994
$ws.Cells.Item(712, 4).Value = 0.005033039709988454
$ws.Cells.Item(713, 4).Value = 0.01155752408321221
$ws.Cells.Item(714, 4).Value = 0.0001043884044200006
$ws.Cells.Item(715, 4).Value = 0.01111580117631435
$ws.Cells.Item(716, 4).Value = 0.0100739749302851
$ws.Cells.Item(717, 4).Value = 0.00587752214872642
$ws.Cells.Item(718, 4).Value = 0.0150976832339956
$ws.Cells.Item(719, 4).Value = 0.03544785938273697
$ws.Cells.Item(720, 4).Value = -0.00225451633946336
$ws.Cells.Item(721, 4).Value = -0.00427172308052727
$ws.Cells.Item(723, 4).Value = 75235437025.98903
$ws.Cells.Item(725, 4).Value = 4359.200631929389
$ws.Cells.Item(727, 4).Value = 199014685312.5529
$ws.Cells.Item(729, 4).Value = 11531.06802154984
$ws.Cells.Item(731, 4).Value = 35707188110.44139
$ws.Cells.Item(733, 4).Value = 2068.902675765529
$ws.Cells.Item(735, 4).Value = 234721873422.9943
$ws.Cells.Item(737, 4).Value = 13599.97069731537
$ws.Cells.Item(739, 4).Value = 58765291864.99162
$ws.Cells.Item(741, 4).Value = 3404.907415436375
$ws.Cells.Item(743, 4).Value = 134000728890.9807
$ws.Cells.Item(745, 4).Value = 7764.108047365764
$ws.Cells.Item(747, 4).Value = 99660392828.35493
$ws.Cells.Item(749, 4).Value = 5774.401858603212
$ws.Cells.Item(751, 4).Value = 1060751703.658653
$ws.Cells.Item(753, 4).Value = 61.46079134639263
$ws.Cells.Item(755, 4).Value = 100721144532.0136
$ws.Cells.Item(757, 4).Value = 5835.862649949605
$ws.Cells.Item(759, 4).Value = 289003532717.5941
$ws.Cells.Item(761, 4).Value = 16745.09290106433
$ws.Cells.Item(763, 4).Value = 5057388683.376679
$ws.Cells.Item(765, 4).Value = 293.0290939477446
$ws.Cells.Item(767, 4).Value = 30649799427.06471
$ws.Cells.Item(769, 4).Value = 1775.873581817784
$ws.Cells.Item(775, 4).Value = 0.1105360742831728
$ws.Cells.Item(776, 4).Value = 0.04575606502754628
$ws.Cells.Item(780, 4).Value = 0.1105980610577872
$ws.Cells.Item(781, 4).Value = 0.04565559689370186
$ws.Cells.Item(785, 4).Value = 0.1140728691162127
$ws.Cells.Item(786, 4).Value = 0.0471756850148608
$ws.Cells.Item(790, 4).Value = 0.1102063733466377
$ws.Cells.Item(791, 4).Value = 0.04547202267905243
$ws.Cells.Item(793, 4).Value = 0.1136472713939958
$ws.Cells.Item(794, 4).Value = 0.04697152250411911
$ws.Cells.Item(798, 4).Value = 0.1136701346791387
$ws.Cells.Item(799, 4).Value = 0.04698189349744572
$ws.Cells.Item(803, 4).Value = 0.1140500058310698
$ws.Cells.Item(804, 4).Value = 0.04716526036804179
$ws.Cells.Item(808, 4).Value = 0.1127981358877047
$ws.Cells.Item(809, 4).Value = 0.04661303399345384
$ws.Cells.Item(813, 4).Value = 0.1122715699427017
$ws.Cells.Item(814, 4).Value = 0.04631724809623655
$ws.Cells.Item(818, 4).Value = 0.1114360955044562
$ws.Cells.Item(819, 4).Value = 0.0459661380219958
$ws.Cells.Item(823, 4).Value = 0.1123899112354324
$ws.Cells.Item(824, 4).Value = 0.04646950677049509
$ws.Cells.Item(828, 4).Value = 0.1136223888967598
$ws.Cells.Item(829, 4).Value = 0.04695954759560618
$ws.Cells.Item(833, 4).Value = 0.1123653716174249
$ws.Cells.Item(834, 4).Value = 0.04645771469160451
